# "added timing module, restructured"
# The substantive change in this workbook revision is a data correction:
# the "Total" figure in B22 loses its trailing footnote marker, going
# from "31 865*" to "31 865" (the row above stays a text value, not a
# number, to preserve the space-separated formatting).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B22").Value = "31 865"

# Reflect the updated view state recorded alongside the edit (the sheet
# was scrolled down and the corrected total cell left selected).
$ws.Range("B22").Select()
